$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "none" / "n/a" fuel row (row 16): a zeroed-out placeholder fuel
# entry, used as a "no fuel selected" option.
$ws.Range("A16").Value = "none"
$ws.Range("B16:P16").Value = 0
$ws.Range("B16:P16").NumberFormat = "0.00"
$ws.Range("Q16").Value = "n/a"

# Match the author's final cursor position/selection.
$ws.Range("F16").Select()
